# Remove the three slides that no longer belong in the deck:
#   "Introduction to Git and GitHub"
#   "Installing Git"
#   "Basic Git Commands"
# The remaining slides (Title, "Working with Branches",
# "Working with GitHub", "Collaboration and Pull Requests",
# "Conclusion") stay in their original relative order.

$p = $ppt.ActivePresentation

$titlesToRemove = @(
    "Introduction to Git and GitHub",
    "Installing Git",
    "Basic Git Commands"
)

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
    if ($titlesToRemove -contains $title) {
        $slide.Delete()
    }
}
